$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2-6 with new values in columns B, C, D (E stays "originario_principal")
$ws.Range("B2").Value = "5202537-66.2022.8.21.0001"
$ws.Range("C2").Value = "5012652-04.2020.8.21.0001"
$ws.Range("D2").Value = "CIV.35317.01"

$ws.Range("B3").Value = "5168190-07.2022.8.21.0001"
$ws.Range("C3").Value = "5016585-82.2020.8.21.0001"
$ws.Range("D3").Value = "CIV.35335.01"

$ws.Range("B4").Value = "5008633-12.2022.8.21.4001"
$ws.Range("C4").Value = "5002299-64.2019.8.21.4001"
$ws.Range("D4").Value = "CIV.13026.01"

$ws.Range("B5").Value = "5008633-12.2022.8.21.4001"
$ws.Range("C5").Value = "5002299-64.2019.8.21.4001"
$ws.Range("D5").Value = "CIV.13026.01"

$ws.Range("B6").Value = "5010628-20.2022.8.21.0005"
$ws.Range("C6").Value = "5003354-73.2020.8.21.0005"
$ws.Range("D6").Value = "CIV.36151.01"

# Delete rows 7 through 11 entirely, shrinking the used range to A1:E6
$ws.Range("A7:E11").EntireRow.Delete()
